# Replace the sample health-tracker data with a new, larger data set and
# re-apply a light banding/border style to the "Name" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rewrite the table contents (header + 6 new people, 7 rows total)
# ---------------------------------------------------------------------
$headers = @("Name", "BloodPressure", "Heartrate", "Sleep", "Date")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value = $headers[$c]
}

$data = @(
    @("Sophia Wilson",       "125/80 mmHg", "72 bpm",  "7 hours",   "11/11/2023"),
    @("Liam Brown",          "130/75 mmHg", "82 bpm",  "5.5 hours", "11/09/2023"),
    @("Olivia Taylor",       "140/85 mmHg", "62 bpm",  "7.5 hours", "11/08/2023"),
    @("Noah Davis",          "70/45 mmHg",  "122 bpm", "9 hours",   "11/02/2023"),
    @("Ella Brown",          "90/60 mmHg",  "67 bpm",  "10 hours",  "11/07/2023"),
    @("Charlotte Anderson",  "135/75 mmHg", "87 bpm",  "8 hours",   "11/08/2023")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $rowData = $data[$r]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = [DateTime]::Parse($rowData[4])
    $ws.Cells.Item($row, 5).NumberFormat = "m/d/yyyy"
}

# ---------------------------------------------------------------------
# 2. Apply a thin green (Accent 6, Lighter 40%) top/bottom border to the
#    "Name" column cells that hold data, plus the blank row right below
#    the table.
# ---------------------------------------------------------------------
$nameRange = $ws.Range("A2:A8")
$nameRange.Borders.Item(8).LineStyle = 1   # xlEdgeTop    -> xlContinuous
$nameRange.Borders.Item(8).Weight = 2      # xlThin
$nameRange.Borders.Item(8).ThemeColor = 9
$nameRange.Borders.Item(8).TintAndShade = 0.39997558519241921

$nameRange.Borders.Item(9).LineStyle = 1   # xlEdgeBottom -> xlContinuous
$nameRange.Borders.Item(9).Weight = 2      # xlThin
$nameRange.Borders.Item(9).ThemeColor = 9
$nameRange.Borders.Item(9).TintAndShade = 0.39997558519241921

# ---------------------------------------------------------------------
# 3. Column / view cosmetics
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.33203125

$ws.Range("E8").Select()
$excel.ActiveWindow.Zoom = 131
